$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 <= content of old row 8
$ws.Range("A7").Value2 = 131085805
$ws.Range("AB7").Value2 = "12:01"
$ws.Range("AW7").Value2 = "Kim Hultgren"
$ws.Range("AX7").Value2 = "Kim Hultgren"
$ws.Range("Q7").Value2 = 585215
$ws.Range("R7").Value2 = 7060513
$ws.Range("S7").Value2 = 10
$ws.Range("Z7").Value2 = "12:01"

# Row 8 <= content of old row 7
$ws.Range("A8").Value2 = 131092646
$ws.Range("AW8").Value2 = "Daniel Rutschman"
$ws.Range("AX8").Value2 = "Daniel Rutschman"
$ws.Range("Q8").Value2 = 585082
$ws.Range("R8").Value2 = 7060264
$ws.Range("S8").Value2 = 15
$ws.Range("AB8").ClearContents()
$ws.Range("Z8").ClearContents()

# Row 18 <= content of old row 20
$ws.Range("A18").Value2 = 131085737
$ws.Range("AB18").Value2 = "11:58"
$ws.Range("AW18").Value2 = "Daniel Rutschman"
$ws.Range("AX18").Value2 = "Daniel Rutschman"
$ws.Range("B18").Value2 = 79244
$ws.Range("E18").Value2 = 6425
$ws.Range("F18").Value2 = "Garnlav"
$ws.Range("G18").Value2 = "Alectoria sarmentosa"
$ws.Range("H18").Value2 = "(Ach.) Ach."
$ws.Range("Q18").Value2 = 585170
$ws.Range("R18").Value2 = 7060469
$ws.Range("S18").Value2 = 15
$ws.Range("Z18").Value2 = "11:58"

# Row 20 <= content of old row 22
$ws.Range("A20").Value2 = 131086957
$ws.Range("AB20").Value2 = "12:21"
$ws.Range("AC20").Value2 = "Ringhack på gran"
$ws.Range("AW20").Value2 = "Kim Hultgren"
$ws.Range("AX20").Value2 = "Kim Hultgren"
$ws.Range("B20").Value2 = 57884
$ws.Range("E20").Value2 = 100109
$ws.Range("F20").Value2 = "Tretåig hackspett"
$ws.Range("G20").Value2 = "Picoides tridactylus"
$ws.Range("H20").Value2 = "(Linnaeus, 1758)"
$ws.Range("M20").Value2 = "färska spår"
$ws.Range("Q20").Value2 = 585162
$ws.Range("R20").Value2 = 7060573
$ws.Range("S20").Value2 = 10
$ws.Range("Z20").Value2 = "12:21"

# Row 21 <= content of old row 18
$ws.Range("A21").Value2 = 131092560
$ws.Range("AB21").Value2 = "15:17"
$ws.Range("AW21").Value2 = "Kim Hultgren"
$ws.Range("AX21").Value2 = "Kim Hultgren"
$ws.Range("B21").Value2 = 91805
$ws.Range("E21").Value2 = 1108
$ws.Range("F21").Value2 = "Harticka"
$ws.Range("G21").Value2 = "Pelloporus leporinus"
$ws.Range("H21").Value2 = "(Fr.) Krieglst."
$ws.Range("Q21").Value2 = 585129
$ws.Range("R21").Value2 = 7060254
$ws.Range("S21").Value2 = 10
$ws.Range("Z21").Value2 = "15:17"
$ws.Range("AC21").ClearContents()
$ws.Range("M21").ClearContents()

# Row 22 <= content of old row 21
$ws.Range("A22").Value2 = 131092554
$ws.Range("AC22").Value2 = "Äldre ringhack, gran"
$ws.Range("AW22").Value2 = "Daniel Rutschman"
$ws.Range("AX22").Value2 = "Daniel Rutschman"
$ws.Range("M22").Value2 = "äldre spår"
$ws.Range("Q22").Value2 = 585147
$ws.Range("R22").Value2 = 7060312
$ws.Range("S22").Value2 = 15
$ws.Range("AB22").ClearContents()
$ws.Range("Z22").ClearContents()

# Row 28 <= content of old row 29
$ws.Range("A28").Value2 = 131085171
$ws.Range("AW28").Value2 = "Daniel Rutschman"
$ws.Range("AX28").Value2 = "Daniel Rutschman"
$ws.Range("B28").Value2 = 91805
$ws.Range("E28").Value2 = 1108
$ws.Range("F28").Value2 = "Harticka"
$ws.Range("G28").Value2 = "Pelloporus leporinus"
$ws.Range("H28").Value2 = "(Fr.) Krieglst."
$ws.Range("Q28").Value2 = 585222
$ws.Range("R28").Value2 = 7060254
$ws.Range("S28").Value2 = 15
$ws.Range("AB28").ClearContents()
$ws.Range("Z28").ClearContents()

# Row 29 <= content of old row 28
$ws.Range("A29").Value2 = 131085178
$ws.Range("AB29").Value2 = "11:08"
$ws.Range("AW29").Value2 = "Kim Hultgren"
$ws.Range("AX29").Value2 = "Kim Hultgren"
$ws.Range("B29").Value2 = 91829
$ws.Range("E29").Value2 = 5432
$ws.Range("F29").Value2 = "Granticka"
$ws.Range("G29").Value2 = "Porodaedalea chrysoloma s.lat."
$ws.Range("H29").ClearContents()
$ws.Range("Q29").Value2 = 585225
$ws.Range("R29").Value2 = 7060258
$ws.Range("S29").Value2 = 10
$ws.Range("Z29").Value2 = "11:08"

# Row 33 <= content of old row 34
$ws.Range("A33").Value2 = 131087388
$ws.Range("B33").Value2 = 79244
$ws.Range("E33").Value2 = 6425
$ws.Range("F33").Value2 = "Garnlav"
$ws.Range("G33").Value2 = "Alectoria sarmentosa"
$ws.Range("H33").Value2 = "(Ach.) Ach."
$ws.Range("Q33").Value2 = 585131
$ws.Range("R33").Value2 = 7060627

# Row 34 <= content of old row 35
$ws.Range("A34").Value2 = 131092590
$ws.Range("AB34").Value2 = "15:20"
$ws.Range("AW34").Value2 = "Kim Hultgren"
$ws.Range("AX34").Value2 = "Kim Hultgren"
$ws.Range("Q34").Value2 = 585145
$ws.Range("R34").Value2 = 7060230
$ws.Range("S34").Value2 = 10
$ws.Range("Z34").Value2 = "15:20"

# Row 35 <= content of old row 33
$ws.Range("A35").Value2 = 131092585
$ws.Range("AW35").Value2 = "Daniel Rutschman"
$ws.Range("AX35").Value2 = "Daniel Rutschman"
$ws.Range("B35").Value2 = 91805
$ws.Range("E35").Value2 = 1108
$ws.Range("F35").Value2 = "Harticka"
$ws.Range("G35").Value2 = "Pelloporus leporinus"
$ws.Range("H35").Value2 = "(Fr.) Krieglst."
$ws.Range("Q35").Value2 = 585130
$ws.Range("R35").Value2 = 7060263
$ws.Range("S35").Value2 = 15
$ws.Range("AB35").ClearContents()
$ws.Range("Z35").ClearContents()
